$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.832.15"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.706.31"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9960"
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.15"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9964"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3924"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4060"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.496"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.65"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9955"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08814"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.26"
$ws.Range("E13").Value = "  +10.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.490"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.112"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001363"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.704.48"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.64"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07168"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.292"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9965"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.40"
$ws.Range("E23").Value = "  -2.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.821.22"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.025"
$ws.Range("E25").Value = "  -4.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.332"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.01"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.66"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.945"
$ws.Range("E29").Value = "  +15.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.572"
$ws.Range("E30").Value = "  -7.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "144.93"
$ws.Range("E31").Value = "  +6.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.890.56"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08833"
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.177"
$ws.Range("E34").Value = "  +10.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.064"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03144"
$ws.Range("E36").Value = "  +5.49%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.231"
$ws.Range("E37").Value = "  -8.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2825"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.94"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8447"
$ws.Range("E40").Value = "  +9.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09221"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.16"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.478"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.63"
$ws.Range("E44").Value = "  +9.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.704"
$ws.Range("E45").Value = "  +4.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7481"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.278"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.397"
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9961"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.65"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08272"
$ws.Range("E51").Value = "  +3.59%  "
